# Wraps a unique phrase already present in the document with
# <w:proofErr w:type="gramStart"/> ... <w:proofErr w:type="gramEnd"/>
# (mirroring what Word's grammar checker stamps onto a flagged phrase),
# without altering the visible text.
#
# InsertXML on a Range replaces that range's own text correctly, but this
# host's Range position bookkeeping places the *newly inserted* nodes at
# the end of the enclosing paragraph rather than in place. To work around
# that we temporarily split the paragraph right after the target phrase
# (so "end of paragraph" and "end of our range" are the same point),
# perform the InsertXML, then delete the now-redundant paragraph mark to
# re-join the paragraph.
function Wrap-GrammarPhrase {
    param(
        [string]$Phrase
    )

    $d = $word.ActiveDocument

    $c = $d.Content
    $found = $c.Find.Execute($Phrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Phrase not found: $Phrase"
    }

    $start = $c.Start
    $end = $c.End

    # Split the paragraph right after the phrase so the range end lines
    # up with a paragraph end (works around the InsertXML placement quirk).
    $splitPoint = $d.Range($end, $end)
    $splitPoint.InsertParagraphAfter()

    $target = $d.Range($start, $end)
    $escaped = $Phrase.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>' + $escaped + '</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($xml)

    # Remove the paragraph mark introduced above to re-merge the paragraph.
    $mark = $d.Range($end, $end + 1)
    $mark.Delete()
}

Wrap-GrammarPhrase "compare and contrast"
Wrap-GrammarPhrase "whether or not"
Wrap-GrammarPhrase "justified, but"
